# Apply "break out stock.yaml completed" changes:
#  - three_line sheet: append row 118
#  - two_line sheet: append rows 16, 17, 18
#  - ph_pl_breakout_line sheet: append rows 498, 499, 500

$wb = $excel.ActiveWorkbook

$dateFmt = "YYYY-MM-DD HH:MM:SS"

function Set-DateCell($ws, $row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = $dateFmt
    $cell.Value2 = $value
}

function Set-Cell($ws, $row, $col, $value) {
    $ws.Cells.Item($row, $col).Value2 = $value
}

# ---------------------------------------------------------------
# Sheet: three_line  -> add row 118 (columns A..L)
# ---------------------------------------------------------------
$wsThree = $wb.Worksheets.Item("three_line")

$r = 118
Set-DateCell $wsThree $r 1  45450.59375          # A118 detected_date
Set-Cell     $wsThree $r 2  "10-06-2024 10:15:00" # B118 breakout_date
Set-Cell     $wsThree $r 3  "hour"                 # C118 Time_Frame
Set-Cell     $wsThree $r 4  "TRF.NS"                # D118 stockname
Set-DateCell $wsThree $r 5  45446.51041666666     # E118 date1
Set-Cell     $wsThree $r 6  580.1500244140625      # F118 value1
Set-DateCell $wsThree $r 7  45446.55208333334     # G118 date2
Set-Cell     $wsThree $r 8  580.1500244140625      # H118 value2
Set-DateCell $wsThree $r 9  45446.63541666666     # I118 date3
Set-Cell     $wsThree $r 10 580.1500244140625      # J118 value3
Set-Cell     $wsThree $r 11 "High"                 # K118 buyORsell
Set-Cell     $wsThree $r 12 "10/06/2024 05:48:26"  # L118 Date Time

# ---------------------------------------------------------------
# Sheet: two_line -> add rows 16, 17, 18 (columns A..J)
# ---------------------------------------------------------------
$wsTwo = $wb.Worksheets.Item("two_line")

$r = 16
Set-DateCell $wsTwo $r 1  45448.38541666666
Set-Cell     $wsTwo $r 2  "10-06-2024 10:15:00"
Set-Cell     $wsTwo $r 3  "hour"
Set-Cell     $wsTwo $r 4  "TRF.NS"
Set-DateCell $wsTwo $r 5  45446.46875
Set-Cell     $wsTwo $r 6  580.1500244140625
Set-DateCell $wsTwo $r 7  45446.51041666666
Set-Cell     $wsTwo $r 8  580.1500244140625
Set-Cell     $wsTwo $r 9  "High"
Set-Cell     $wsTwo $r 10 "10/06/2024 05:48:26"

$r = 17
Set-DateCell $wsTwo $r 1  45448.42708333334
Set-Cell     $wsTwo $r 2  "10-06-2024 10:15:00"
Set-Cell     $wsTwo $r 3  "hour"
Set-Cell     $wsTwo $r 4  "TRF.NS"
Set-DateCell $wsTwo $r 5  45446.46875
Set-Cell     $wsTwo $r 6  580.1500244140625
Set-DateCell $wsTwo $r 7  45446.55208333334
Set-Cell     $wsTwo $r 8  580.1500244140625
Set-Cell     $wsTwo $r 9  "High"
Set-Cell     $wsTwo $r 10 "10/06/2024 05:48:26"

$r = 18
Set-DateCell $wsTwo $r 1  45448.42708333334
Set-Cell     $wsTwo $r 2  "10-06-2024 10:15:00"
Set-Cell     $wsTwo $r 3  "hour"
Set-Cell     $wsTwo $r 4  "TRF.NS"
Set-DateCell $wsTwo $r 5  45446.51041666666
Set-Cell     $wsTwo $r 6  580.1500244140625
Set-DateCell $wsTwo $r 7  45446.55208333334
Set-Cell     $wsTwo $r 8  580.1500244140625
Set-Cell     $wsTwo $r 9  "High"
Set-Cell     $wsTwo $r 10 "10/06/2024 05:48:26"

# ---------------------------------------------------------------
# Sheet: ph_pl_breakout_line -> add rows 498, 499, 500 (columns A..L)
# ---------------------------------------------------------------
$wsPh = $wb.Worksheets.Item("ph_pl_breakout_line")

$r = 498
Set-Cell     $wsPh $r 1  "FACT.NS"
Set-DateCell $wsPh $r 2  45436.42708333334
Set-Cell     $wsPh $r 3  747.7999877929688
Set-Cell     $wsPh $r 4  732.4500122070312
Set-Cell     $wsPh $r 5  736.7000122070312
Set-Cell     $wsPh $r 6  "High"
Set-Cell     $wsPh $r 7  747.7999877929688
Set-Cell     $wsPh $r 8  "hour"
Set-Cell     $wsPh $r 9  "10-06-2024 11:15:00"
Set-Cell     $wsPh $r 10 759.5499877929688
Set-Cell     $wsPh $r 11 747.5999755859375
Set-Cell     $wsPh $r 12 "10/06/2024 05:48:26"

$r = 499
Set-Cell     $wsPh $r 1  "TV18BRDCST.NS"
Set-DateCell $wsPh $r 2  45446.38541666666
Set-Cell     $wsPh $r 3  43.59999847412109
Set-Cell     $wsPh $r 4  43
Set-Cell     $wsPh $r 5  43.20000076293945
Set-Cell     $wsPh $r 6  "High"
Set-Cell     $wsPh $r 7  43.59999847412109
Set-Cell     $wsPh $r 8  "hour"
Set-Cell     $wsPh $r 9  "10-06-2024 11:15:00"
Set-Cell     $wsPh $r 10 43.59999847412109
Set-Cell     $wsPh $r 11 43.52000045776367
Set-Cell     $wsPh $r 12 "10/06/2024 05:48:26"

$r = 500
Set-Cell     $wsPh $r 1  "AUROPHARMA.NS"
Set-DateCell $wsPh $r 2  45436.59375
Set-Cell     $wsPh $r 3  1245.699951171875
Set-Cell     $wsPh $r 4  1234
Set-Cell     $wsPh $r 5  1235.150024414062
Set-Cell     $wsPh $r 6  "High"
Set-Cell     $wsPh $r 7  1245.699951171875
Set-Cell     $wsPh $r 8  "hour"
Set-Cell     $wsPh $r 9  "10-06-2024 11:15:00"
Set-Cell     $wsPh $r 10 1245.900024414062
Set-Cell     $wsPh $r 11 1244.75
Set-Cell     $wsPh $r 12 "10/06/2024 05:48:26"
